# Apply scheduled-runner profit recalculations to the Odin_Profits sheets.
# Each row below corresponds to a Leve whose average-price inputs (and the
# downstream Leve price / profit columns) were refreshed by the pricing bot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2430.6223
$ws.Range("J17").Value = 2435.6584
$ws.Range("L17").Value = 7306.975199999999
$ws.Range("N17").Value = -7642.975199999999

$ws.Range("H40").Value = 9663.666999999999
$ws.Range("J40").Value = 9995
$ws.Range("L40").Value = 9995
$ws.Range("N40").Value = -10345

$ws.Range("H132").Value = 279782.94
$ws.Range("I132").Value = 337092.47
$ws.Range("K132").Value = 1011277.41
$ws.Range("M132").Value = -1008747.41

$ws.Range("H137").Value = 6653.0527
$ws.Range("I137").Value = 9320.333000000001
$ws.Range("J137").Value = 5422
$ws.Range("K137").Value = 27960.999
$ws.Range("L137").Value = 16266
$ws.Range("M137").Value = -25410.999
$ws.Range("N137").Value = -21366

$ws.Range("H138").Value = 4031.4058
$ws.Range("J138").Value = 4669.3555
$ws.Range("L138").Value = 14008.0665
$ws.Range("N138").Value = -24288.0665

$ws.Range("H139").Value = 179998.67
$ws.Range("J139").Value = 179998.67
$ws.Range("L139").Value = 179998.67
$ws.Range("N139").Value = -190278.67

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7315.16
$ws.Range("I32").Value = 3133.4604
$ws.Range("J32").Value = 20557.209
$ws.Range("K32").Value = 3133.4604
$ws.Range("L32").Value = 20557.209
$ws.Range("M32").Value = -2846.4604
$ws.Range("N32").Value = -21131.209

$ws.Range("H74").Value = 4853.6553
$ws.Range("I74").Value = 4873.5557
$ws.Range("K74").Value = 4873.5557
$ws.Range("M74").Value = -3999.5557

$ws.Range("H77").Value = 4853.6553
$ws.Range("I77").Value = 4873.5557
$ws.Range("K77").Value = 24367.7785
$ws.Range("M77").Value = -19999.7785

$ws.Range("H86").Value = 131400
$ws.Range("J86").Value = 131400
$ws.Range("L86").Value = 131400
$ws.Range("N86").Value = -133772

$ws.Range("H89").Value = 131400
$ws.Range("J89").Value = 131400
$ws.Range("L89").Value = 394200
$ws.Range("N89").Value = -406056

$ws.Range("H110").Value = 6437.75
$ws.Range("I110").Value = 4324.125
$ws.Range("J110").Value = 7283.2
$ws.Range("K110").Value = 4324.125
$ws.Range("L110").Value = 7283.2
$ws.Range("M110").Value = -2279.125
$ws.Range("N110").Value = -11373.2

$ws.Range("H123").Value = 70000
$ws.Range("J123").Value = 70000
$ws.Range("L123").Value = 70000
$ws.Range("N123").Value = -79800

$ws.Range("H139").Value = 103984.5
$ws.Range("J139").Value = 103984.5
$ws.Range("L139").Value = 103984.5
$ws.Range("N139").Value = -114264.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 8091.092
$ws.Range("I99").Value = 7173.517
$ws.Range("J99").Value = 9926.241
$ws.Range("K99").Value = 7173.517
$ws.Range("L99").Value = 9926.241
$ws.Range("M99").Value = -5675.517
$ws.Range("N99").Value = -12922.241

$ws.Range("H105").Value = 4407
$ws.Range("J105").Value = 3499
$ws.Range("L105").Value = 3499
$ws.Range("N105").Value = -6993

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6155.048
$ws.Range("I31").Value = 1971.8889
$ws.Range("J31").Value = 9292.416999999999
$ws.Range("K31").Value = 1971.8889
$ws.Range("L31").Value = 9292.416999999999
$ws.Range("M31").Value = -1676.8889
$ws.Range("N31").Value = -9882.416999999999

$ws.Range("H34").Value = 6155.048
$ws.Range("I34").Value = 1971.8889
$ws.Range("J34").Value = 9292.416999999999
$ws.Range("K34").Value = 1971.8889
$ws.Range("L34").Value = 9292.416999999999
$ws.Range("M34").Value = -1769.8889
$ws.Range("N34").Value = -9696.416999999999

$ws.Range("H132").Value = 7648.273
$ws.Range("I132").Value = 5719.346
$ws.Range("J132").Value = 10434.5
$ws.Range("K132").Value = 17158.038
$ws.Range("L132").Value = 31303.5
$ws.Range("M132").Value = -14628.038
$ws.Range("N132").Value = -36363.5

$ws.Range("H138").Value = 85250
$ws.Range("I138").Value = 100000
$ws.Range("J138").Value = 70500
$ws.Range("K138").Value = 100000
$ws.Range("L138").Value = 70500
$ws.Range("M138").Value = -94860
$ws.Range("N138").Value = -80780

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 113779.27
$ws.Range("J37").Value = 113779.27
$ws.Range("L37").Value = 341337.81
$ws.Range("N37").Value = -341561.81

$ws.Range("H64").Value = 20470.385
$ws.Range("I64").Value = 5000
$ws.Range("J64").Value = 21759.584
$ws.Range("K64").Value = 15000
$ws.Range("L64").Value = 65278.75199999999
$ws.Range("M64").Value = -14730
$ws.Range("N64").Value = -65818.75199999999

$ws.Range("H67").Value = 20470.385
$ws.Range("I67").Value = 5000
$ws.Range("J67").Value = 21759.584
$ws.Range("K67").Value = 15000
$ws.Range("L67").Value = 65278.75199999999
$ws.Range("M67").Value = -14064
$ws.Range("N67").Value = -67150.75199999999

$ws.Range("H80").Value = 9100.200000000001
$ws.Range("I80").Value = 7751
$ws.Range("J80").Value = 9999.666999999999
$ws.Range("K80").Value = 23253
$ws.Range("L80").Value = 29999.001
$ws.Range("M80").Value = -22317
$ws.Range("N80").Value = -31871.001

$ws.Range("H83").Value = 9100.200000000001
$ws.Range("I83").Value = 7751
$ws.Range("J83").Value = 9999.666999999999
$ws.Range("K83").Value = 69759
$ws.Range("L83").Value = 89997.003
$ws.Range("M83").Value = -65079
$ws.Range("N83").Value = -99357.003

$ws.Range("H114").Value = 1316.5883
$ws.Range("I114").Value = 979.0909
$ws.Range("J114").Value = 1935.3334
$ws.Range("K114").Value = 2937.2727
$ws.Range("L114").Value = 5806.0002
$ws.Range("M114").Value = 316.7273
$ws.Range("N114").Value = -12314.0002

$ws.Range("H117").Value = 4214.5
$ws.Range("I117").Value = 1424.4
$ws.Range("K117").Value = 4273.200000000001
$ws.Range("M117").Value = -831.2000000000007

$ws.Range("H136").Value = 55561556
$ws.Range("J136").Value = 10000
$ws.Range("L136").Value = 30000
$ws.Range("N136").Value = -40200

$ws.Range("H139").Value = 2912.6155
$ws.Range("I139").Value = 2169
$ws.Range("J139").Value = 4585.75
$ws.Range("K139").Value = 6507
$ws.Range("L139").Value = 13757.25
$ws.Range("M139").Value = -1367
$ws.Range("N139").Value = -24037.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 31256330
$ws.Range("I80").Value = 125003370
$ws.Range("J80").Value = 7318.7085
$ws.Range("K80").Value = 125003370
$ws.Range("L80").Value = 7318.7085
$ws.Range("M80").Value = -125002372
$ws.Range("N80").Value = -9314.708500000001

$ws.Range("H83").Value = 31256330
$ws.Range("I83").Value = 125003370
$ws.Range("J83").Value = 7318.7085
$ws.Range("K83").Value = 625016850
$ws.Range("L83").Value = 36593.5425
$ws.Range("M83").Value = -625011858
$ws.Range("N83").Value = -46577.5425

$ws.Range("H102").Value = 4321.8887
$ws.Range("I102").Value = 3191
$ws.Range("J102").Value = 5931.231
$ws.Range("K102").Value = 3191
$ws.Range("L102").Value = 5931.231
$ws.Range("M102").Value = -1569
$ws.Range("N102").Value = -9175.231

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 35715676
$ws.Range("I46").Value = 963.4286
$ws.Range("K46").Value = 963.4286
$ws.Range("M46").Value = -775.4286

$ws.Range("H100").Value = 3646.8696
$ws.Range("I100").Value = 5158.4
$ws.Range("J100").Value = 2484.1538
$ws.Range("K100").Value = 5158.4
$ws.Range("L100").Value = 2484.1538
$ws.Range("M100").Value = -4617.4
$ws.Range("N100").Value = -3566.1538

$ws.Range("H139").Value = 69249.75
$ws.Range("J139").Value = 69249.75
$ws.Range("L139").Value = 69249.75
$ws.Range("N139").Value = -79529.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 2012800
$ws.Range("I9").Value = 2012800
$ws.Range("K9").Value = 2012800
$ws.Range("M9").Value = -2012660
